$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1577.3334
$ws.Range("I28").Value = 1737.4445
$ws.Range("K28").Value = 1737.4445
$ws.Range("M28").Value = -1252.4445

$ws.Range("H33").Value = 367.15384
$ws.Range("I33").Value = 215.65517
$ws.Range("K33").Value = 215.65517
$ws.Range("M33").Value = 13.34483

$ws.Range("H38").Value = 4079.5
$ws.Range("I38").Value = 439.33334
$ws.Range("J38").Value = 15000
$ws.Range("K38").Value = 1318.00002
$ws.Range("L38").Value = 45000
$ws.Range("M38").Value = -946.0000199999999
$ws.Range("N38").Value = -45744

$ws.Range("H64").Value = 3763.818
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996

$ws.Range("H67").Value = 3763.818
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216

$ws.Range("H106").Value = 2900403
$ws.Range("I106").Value = 3176155.8
$ws.Range("K106").Value = 3176155.8
$ws.Range("M106").Value = -3175524.8

$ws.Range("H107").Value = 8058.846
$ws.Range("I107").Value = 8475.5
$ws.Range("K107").Value = 8475.5
$ws.Range("M107").Value = -6555.5

$ws.Range("H132").Value = 3917.4348
$ws.Range("I132").Value = 4088.4211
$ws.Range("K132").Value = 12265.2633
$ws.Range("M132").Value = -9735.263300000001

$ws.Range("H138").Value = 2798.6338
$ws.Range("I138").Value = 1940.5385
$ws.Range("J138").Value = 3294.422
$ws.Range("K138").Value = 5821.6155
$ws.Range("L138").Value = 9883.266
$ws.Range("M138").Value = -681.6154999999999
$ws.Range("N138").Value = -20163.266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8014387
$ws.Range("I32").Value = 3473294.5
$ws.Range("J32").Value = 62507500
$ws.Range("K32").Value = 3473294.5
$ws.Range("L32").Value = 62507500
$ws.Range("M32").Value = -3473007.5
$ws.Range("N32").Value = -62508074

$ws.Range("H61").Value = 2494.8572
$ws.Range("I61").Value = 2209.2173
$ws.Range("K61").Value = 2209.2173
$ws.Range("M61").Value = -1997.2173

$ws.Range("H74").Value = 1711.7142
$ws.Range("I74").Value = 1413.6666
$ws.Range("K74").Value = 1413.6666
$ws.Range("M74").Value = -539.6666

$ws.Range("H77").Value = 1711.7142
$ws.Range("I77").Value = 1413.6666
$ws.Range("K77").Value = 7068.333000000001
$ws.Range("M77").Value = -2700.333000000001

$ws.Range("H102").Value = 1153.238
$ws.Range("I102").Value = 1012.1667
$ws.Range("K102").Value = 1012.1667
$ws.Range("M102").Value = 609.8333

$ws.Range("H110").Value = 1934.2307
$ws.Range("I110").Value = 1762.0834
$ws.Range("K110").Value = 1762.0834
$ws.Range("M110").Value = 282.9166

$ws.Range("H132").Value = 2670.2727
$ws.Range("I132").Value = 2193.5151
$ws.Range("K132").Value = 6580.5453
$ws.Range("M132").Value = -4050.5453

$ws.Range("H136").Value = 2494.8572
$ws.Range("I136").Value = 2209.2173
$ws.Range("K136").Value = 6627.651899999999
$ws.Range("M136").Value = -4077.651899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1210.5714
$ws.Range("I107").Value = 1076.8334
$ws.Range("K107").Value = 1076.8334
$ws.Range("M107").Value = 843.1666

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1825.9286
$ws.Range("I16").Value = 1893.6666
$ws.Range("K16").Value = 1893.6666
$ws.Range("M16").Value = -1606.6666

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H31").Value = 2950.4062
$ws.Range("I31").Value = 2110.9546
$ws.Range("K31").Value = 2110.9546
$ws.Range("M31").Value = -1815.9546

$ws.Range("H34").Value = 2950.4062
$ws.Range("I34").Value = 2110.9546
$ws.Range("K34").Value = 2110.9546
$ws.Range("M34").Value = -1908.9546

$ws.Range("H99").Value = 1526
$ws.Range("I99").Value = 1526
$ws.Range("K99").Value = 1526
$ws.Range("M99").Value = -28

$ws.Range("H113").Value = 1825.9286
$ws.Range("I113").Value = 1893.6666
$ws.Range("K113").Value = 1893.6666
$ws.Range("M113").Value = 276.3334

$ws.Range("H126").Value = 1526
$ws.Range("I126").Value = 1526
$ws.Range("K126").Value = 4578
$ws.Range("M126").Value = -2108

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1014.64703
$ws.Range("I5").Value = 619.625
$ws.Range("J5").Value = 1365.7778
$ws.Range("K5").Value = 1858.875
$ws.Range("L5").Value = 4097.3334
$ws.Range("M5").Value = -1746.875
$ws.Range("N5").Value = -4321.3334

$ws.Range("H86").Value = 189.85715
$ws.Range("I86").Value = 183
$ws.Range("K86").Value = 549
$ws.Range("M86").Value = 637

$ws.Range("H89").Value = 189.85715
$ws.Range("I89").Value = 183
$ws.Range("K89").Value = 1647
$ws.Range("M89").Value = 4281

$ws.Range("H108").Value = 3000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 9000
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -14760

$ws.Range("H120").Value = 20471.75
$ws.Range("J120").Value = 24800
$ws.Range("L120").Value = 74400
$ws.Range("N120").Value = -84076

$ws.Range("H135").Value = 1014.64703
$ws.Range("I135").Value = 619.625
$ws.Range("J135").Value = 1365.7778
$ws.Range("K135").Value = 5576.625
$ws.Range("L135").Value = 12292.0002
$ws.Range("M135").Value = -3041.625
$ws.Range("N135").Value = -17362.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 703.64703
$ws.Range("I97").Value = 763.5
$ws.Range("J97").Value = 560
$ws.Range("K97").Value = 763.5
$ws.Range("L97").Value = 560
$ws.Range("M97").Value = -267.5
$ws.Range("N97").Value = -1552

$ws.Range("H113").Value = 1790.5
$ws.Range("I113").Value = 1404.3334
$ws.Range("J113").Value = 2287
$ws.Range("K113").Value = 1404.3334
$ws.Range("L113").Value = 2287
$ws.Range("M113").Value = 765.6666
$ws.Range("N113").Value = -6627

$ws.Range("H126").Value = 3500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1537.7727
$ws.Range("I61").Value = 1515.5238
$ws.Range("K61").Value = 1515.5238
$ws.Range("M61").Value = -1313.5238

$ws.Range("H68").Value = 1680.2
$ws.Range("I68").Value = 1569.1818
$ws.Range("J68").Value = 1985.5
$ws.Range("K68").Value = 1569.1818
$ws.Range("L68").Value = 1985.5
$ws.Range("M68").Value = -820.1818000000001
$ws.Range("N68").Value = -3483.5

$ws.Range("H71").Value = 1680.2
$ws.Range("I71").Value = 1569.1818
$ws.Range("J71").Value = 1985.5
$ws.Range("K71").Value = 7845.909000000001
$ws.Range("L71").Value = 9927.5
$ws.Range("M71").Value = -4101.909000000001
$ws.Range("N71").Value = -17415.5

$ws.Range("H82").Value = 4668.375
$ws.Range("I82").Value = 4148.3335
$ws.Range("K82").Value = 4148.3335
$ws.Range("M82").Value = -3787.3335

$ws.Range("H85").Value = 4668.375
$ws.Range("I85").Value = 4148.3335
$ws.Range("K85").Value = 4148.3335
$ws.Range("M85").Value = -2900.3335

$ws.Range("H113").Value = 1537.7727
$ws.Range("I113").Value = 1515.5238
$ws.Range("K113").Value = 1515.5238
$ws.Range("M113").Value = 654.4762000000001

$ws.Range("H136").Value = 4022.6667
$ws.Range("I136").Value = 3687.4285
$ws.Range("J136").Value = 4492
$ws.Range("K136").Value = 11062.2855
$ws.Range("M136").Value = -8512.2855
$ws.Range("N136").Value = -18576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3192.85
$ws.Range("I81").Value = 2021.3334
$ws.Range("K81").Value = 4042.6668
$ws.Range("M81").Value = -2981.6668

$ws.Range("H84").Value = 3192.85
$ws.Range("I84").Value = 2021.3334
$ws.Range("K84").Value = 20213.334
$ws.Range("M84").Value = -14909.334
